$wb = $excel.ActiveWorkbook

# This script applies cached market-data refresh values to the Ifrit_Profits
# workbook (multiple job/craft sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each sheet has columns H..N holding scraped market-board figures; a scheduled
# runner refreshed these numbers. Cells that no longer have a value are cleared
# entirely (matching how the source XML drops the <c> element).

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 31251200
$ws.Range("J112").Value = 33334578
$ws.Range("L112").Value = 100003734
$ws.Range("N112").Value = -100005950
$ws.Range("H138").Value = 3540.5715
$ws.Range("I138").Value = 3365.8948
$ws.Range("J138").Value = 5200
$ws.Range("K138").Value = 10097.6844
$ws.Range("L138").Value = 15600
$ws.Range("M138").Value = -4957.6844
$ws.Range("N138").Value = -25880

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H32").Value = 3752.3132
$ws.Range("I32").Value = 3991.5225
$ws.Range("J32").Value = 2750.625
$ws.Range("K32").Value = 3991.5225
$ws.Range("L32").Value = 2750.625
$ws.Range("M32").Value = -3704.5225
$ws.Range("N32").Value = -3324.625
$ws.Range("H37").Value = 10224.091
$ws.Range("J37").Value = 11196.5
$ws.Range("L37").Value = 11196.5
$ws.Range("N37").Value = -11742.5
$ws.Range("H44").Value = 19750
$ws.Range("J44").Value = 19750
$ws.Range("L44").Value = 19750
$ws.Range("N44").Value = -20726
$ws.Range("H55").Value = 27507.625
$ws.Range("I55").Value = 20000
$ws.Range("J55").Value = 28580.143
$ws.Range("K55").Value = 20000
$ws.Range("L55").Value = 28580.143
$ws.Range("M55").Value = -19685
$ws.Range("N55").Value = -29210.143
$ws.Range("H61").Value = 2446.4285
$ws.Range("I61").Value = 1600
$ws.Range("J61").Value = 3575
$ws.Range("K61").Value = 1600
$ws.Range("L61").Value = 3575
$ws.Range("M61").Value = -1388
$ws.Range("N61").Value = -3999
$ws.Range("H63").Value = 6999.8
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 6999.8
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 6999.8
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -8371.799999999999
$ws.Range("H66").Value = 6999.8
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 6999.8
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 34999
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -41863
$ws.Range("H80").Value = 33032.5
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 33032.5
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 33032.5
$ws.Range("N80").Value = -35028.5
$ws.Range("H83").Value = 33032.5
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 33032.5
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 99097.5
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -109081.5
$ws.Range("H136").Value = 2446.4285
$ws.Range("I136").Value = 1600
$ws.Range("J136").Value = 3575
$ws.Range("K136").Value = 4800
$ws.Range("L136").Value = 10725
$ws.Range("M136").Value = -2250
$ws.Range("N136").Value = -15825

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 12000
$ws.Range("J19").Value = 12000
$ws.Range("L19").Value = 12000
$ws.Range("N19").Value = -12346
$ws.Range("H140").Value = 36599.8
$ws.Range("J140").Value = 36599.8
$ws.Range("L140").Value = 36599.8
$ws.Range("N140").Value = -46959.8

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 260.48837
$ws.Range("I22").Value = 230.17949
$ws.Range("K22").Value = 230.17949
$ws.Range("M22").Value = 119.82051
$ws.Range("H31").Value = 3525
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 3525
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 3525
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -4115
$ws.Range("H34").Value = 3525
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 3525
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 3525
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -3929
$ws.Range("H62").Value = 3122.1
$ws.Range("I62").Value = 3284.8333
$ws.Range("J62").Value = 2878
$ws.Range("K62").Value = 3284.8333
$ws.Range("L62").Value = 2878
$ws.Range("M62").Value = -2660.8333
$ws.Range("N62").Value = -4126
$ws.Range("H65").Value = 3122.1
$ws.Range("I65").Value = 3284.8333
$ws.Range("J65").Value = 2878
$ws.Range("K65").Value = 16424.1665
$ws.Range("L65").Value = 14390
$ws.Range("M65").Value = -13304.1665
$ws.Range("N65").Value = -20630
$ws.Range("H99").Value = 1097.2
$ws.Range("I99").Value = 1085.6666
$ws.Range("J99").Value = 1201
$ws.Range("K99").Value = 1085.6666
$ws.Range("L99").Value = 1201
$ws.Range("M99").Value = 412.3334
$ws.Range("N99").Value = -4197
$ws.Range("H105").Value = 1167.7778
$ws.Range("J105").Value = 1005
$ws.Range("L105").Value = 1005
$ws.Range("N105").Value = -4499
$ws.Range("H126").Value = 1097.2
$ws.Range("I126").Value = 1085.6666
$ws.Range("J126").Value = 1201
$ws.Range("K126").Value = 3256.9998
$ws.Range("L126").Value = 3603
$ws.Range("M126").Value = -786.9998000000001
$ws.Range("N126").Value = -8543
$ws.Range("H132").Value = 3975.7368
$ws.Range("I132").Value = 3752.6667
$ws.Range("J132").Value = 4358.143
$ws.Range("K132").Value = 11258.0001
$ws.Range("L132").Value = 13074.429
$ws.Range("M132").Value = -8728.000100000001
$ws.Range("N132").Value = -18134.429
$ws.Range("H134").Value = 1735.92
$ws.Range("I134").Value = 1960.1111
$ws.Range("J134").Value = 1159.4286
$ws.Range("K134").Value = 5880.3333
$ws.Range("L134").Value = 3478.2858
$ws.Range("M134").Value = -3345.3333
$ws.Range("N134").Value = -8548.2858

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 3890
$ws.Range("J59").Value = 3890
$ws.Range("L59").Value = 11670
$ws.Range("N59").Value = -12750
$ws.Range("H121").Value = 35715244
$ws.Range("I121").Value = 336.4
$ws.Range("J121").Value = 55556860
$ws.Range("K121").Value = 1009.2
$ws.Range("L121").Value = 166670580
$ws.Range("M121").Value = 300.8000000000001
$ws.Range("N121").Value = -166673200
$ws.Range("H131").Value = 2635430.8
$ws.Range("I131").Value = 4910.385
$ws.Range("J131").Value = 4003301.2
$ws.Range("K131").Value = 14731.155
$ws.Range("L131").Value = 12009903.6
$ws.Range("M131").Value = -9691.155000000001
$ws.Range("N131").Value = -12019983.6

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 541
$ws.Range("I13").Value = 52.5
$ws.Range("J13").Value = 866.6667
$ws.Range("K13").Value = 52.5
$ws.Range("L13").Value = 866.6667
$ws.Range("M13").Value = 86.5
$ws.Range("N13").Value = -1144.6667
$ws.Range("H102").Value = 1495.8667
$ws.Range("I102").Value = 1130.3636
$ws.Range("K102").Value = 1130.3636
$ws.Range("M102").Value = 491.6364000000001

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1025.75
$ws.Range("I61").Value = 967.6667
$ws.Range("J61").Value = 1200
$ws.Range("K61").Value = 967.6667
$ws.Range("L61").Value = 1200
$ws.Range("M61").Value = -765.6667
$ws.Range("N61").Value = -1604
$ws.Range("H113").Value = 1025.75
$ws.Range("I113").Value = 967.6667
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 967.6667
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 1202.3333
$ws.Range("N113").Value = -5540
$ws.Range("H122").Value = 4952.091
$ws.Range("I122").Value = 6076.1787
$ws.Range("J122").Value = 2984.9375
$ws.Range("K122").Value = 18228.5361
$ws.Range("L122").Value = 8954.8125
$ws.Range("M122").Value = -15778.5361
$ws.Range("N122").Value = -13854.8125

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 15399.6
$ws.Range("I62").Value = 2999
$ws.Range("J62").Value = 18499.75
$ws.Range("K62").Value = 2999
$ws.Range("L62").Value = 18499.75
$ws.Range("M62").Value = -2375
$ws.Range("N62").Value = -19747.75
$ws.Range("H65").Value = 15399.6
$ws.Range("I65").Value = 2999
$ws.Range("J65").Value = 18499.75
$ws.Range("K65").Value = 14995
$ws.Range("L65").Value = 92498.75
$ws.Range("M65").Value = -11875
$ws.Range("N65").Value = -98738.75
$ws.Range("H100").Value = 287.61905
$ws.Range("I100").Value = 206
$ws.Range("J100").Value = 491.66666
$ws.Range("K100").Value = 412
$ws.Range("L100").Value = 983.33332
$ws.Range("M100").Value = 129
$ws.Range("N100").Value = -2065.33332
$ws.Range("H122").Value = 1652
$ws.Range("I122").Value = 1084.2
$ws.Range("J122").Value = 2219.8
$ws.Range("K122").Value = 3252.6
$ws.Range("L122").Value = 6659.400000000001
$ws.Range("M122").Value = -802.6000000000004
$ws.Range("N122").Value = -11559.4
$ws.Range("H126").Value = 5404
$ws.Range("I126").Value = 5404
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 16212
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -13742
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 10410.723
$ws.Range("I132").Value = 11087.1875
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 33261.5625
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -30731.5625
$ws.Range("N132").Value = -20057
